$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The rows' variable data (columns D, L, M, N, O, P, Q, R, S, T) is
# cyclically rotated among rows 2, 4, 5, 6, 7, 8, 9, 10 (rows 3 and 11
# are unaffected). We capture the original values first, then write
# the rotated values back so we don't clobber data before reading it.

$rows = @(2, 4, 5, 6, 7, 8, 9, 10)

# Mapping: destination row -> source row (source row's original data
# becomes the destination row's new data)
$mapping = @{
    2  = 6
    4  = 2
    5  = 4
    6  = 9
    7  = 10
    8  = 5
    9  = 7
    10 = 8
}

$cols = @("D", "L", "M", "N", "O", "P", "Q", "R", "S", "T")

# Snapshot original values for the affected rows/columns
# (use Value2 for reads; Value's getter is not usable in this runtime)
$orig = @{}
foreach ($r in $rows) {
    $orig[$r] = @{}
    foreach ($c in $cols) {
        $orig[$r][$c] = $ws.Range("$c$r").Value2
    }
}

# Apply the rotated values
foreach ($r in $rows) {
    $src = $mapping[$r]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value = $orig[$src][$c]
    }
}
